# edit.ps1
# Applies the Foxlocket.docx revision described in the commit "Again: still did not
# changed to dynamic addresses." to the document currently open as $word.ActiveDocument.
#
# Summary of the edit (see task diff):
#   1. Wrap the English abbreviation "tp" in <w:proofErr spellStart/spellEnd> markers.
#   2. Wrap "SetSensitivity" (in a table cell) the same way.
#   3. Wrap the first occurrence of "CarrierDetect" the same way.
#   4. Merge several runs of the "Сменил множитель..." paragraph into one run, and
#      also proofErr-wrap the (now second) occurrence of "CarrierDetect" in it, splitting
#      off the trailing "." into its own (still en-US) run.
#   5. Rewrite the final part of the document: the old 3 trailing paragraphs
#      ("...А это неверно." / "Попробуем учесть..." / "Нужно") are replaced by a new
#      sequence of 10 paragraphs (new diary entries), and the _GoBack bookmark moves from
#      the first of those paragraphs to a brand-new, otherwise empty, final paragraph.
#
# Implementation note: this COM-interop runtime's Range.InsertXML always behaves as a
# paragraph-level (block) insertion - it cannot splice a bare, non-text element such as
# <w:proofErr/> into the middle of an existing paragraph without Word rebuilding that
# paragraph. So every change below is applied by selecting the *whole* paragraph(s) that
# need to change and replacing them wholesale with freshly authored OOXML that already
# contains the desired <w:proofErr/> markers / merged runs / new paragraphs. Paragraphs
# that are not mentioned in the diff are left completely untouched.

$d = $word.ActiveDocument

function Get-ParagraphByText($doc, [string]$needle) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

function Replace-ParagraphXml($rng, [string]$innerXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# --- Change 1: "В дальнейшем длительность пакета обозначается tp." ---------------------
$target1 = @'
<w:p w14:paraId="10676099" w14:textId="6DB64A05" w:rsidR="001C2B35" w:rsidRDefault="001C2B35" w:rsidP="005E663D"><w:r><w:t xml:space="preserve">В дальнейшем длительность пакета обозначается </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>tp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="001C2B35"><w:t>.</w:t></w:r><w:r w:rsidR="00612CEB"><w:t xml:space="preserve"> Другой – передающий – медальон обозначается МА; медальон, о коем идет речь, обозначается МБ.</w:t></w:r></w:p>
'@
$p1 = Get-ParagraphByText $d "В дальнейшем длительность пакета обозначается"
Replace-ParagraphXml $p1.Range $target1

# --- Change 2: "SetSensitivity" table-cell paragraph --------------------------------
$target2 = @'
<w:p w14:paraId="1E4BEECA" w14:textId="0E9A5E68" w:rsidR="005E663D" w:rsidRPr="004B6F8A" w:rsidRDefault="004B6F8A" w:rsidP="005E663D"><w:pPr><w:ind w:firstLine="0"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>SetSensitivity</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$p2 = Get-ParagraphByText $d "SetSensitivity"
Replace-ParagraphXml $p2.Range $target2

# --- Change 3: first "CarrierDetect" occurrence --------------------------------------
$target3 = @'
<w:p w14:paraId="08A56379" w14:textId="788B65EA" w:rsidR="00F25CA0" w:rsidRPr="0001799F" w:rsidRDefault="00A93646" w:rsidP="00F25CA0"><w:r><w:t xml:space="preserve">Еще нужно уменьшить множитель, поскольку 27 мс – это многовато. В 200 мс 10 медальонов не влезут. Тем более что </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>CarrierDetect</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00F25CA0"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>говорит, что длина передачи 6 мс.</w:t></w:r><w:r w:rsidR="00F25CA0" w:rsidRPr="00F25CA0"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00F25CA0"><w:t xml:space="preserve">Но первым делом – отключить </w:t></w:r><w:r w:rsidR="00F25CA0"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>UART</w:t></w:r><w:r w:rsidR="00F25CA0" w:rsidRPr="0001799F"><w:t>.</w:t></w:r></w:p>
'@
$p3 = Get-ParagraphByText $d "Еще нужно уменьшить множитель"
Replace-ParagraphXml $p3.Range $target3

# --- Change 4: "Сменил множитель с 8 ..." paragraph (merge runs + proofErr) ----------
$target4 = @'
<w:p w14:paraId="7C1E264A" w14:textId="13D116F8" w:rsidR="00E3559F" w:rsidRDefault="00E3559F" w:rsidP="00D203E1"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Сменил множитель с 8 на 4. Расстояние уменьшилось до 13.4 мс. Принятый пакет поступает через 7.4 мс после прерывания на передачу. Длительность пакета примерно 6.8 мс – измерено по </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>CarrierDetect</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@
$p4 = Get-ParagraphByText $d "Сменил множитель с 8"
Replace-ParagraphXml $p4.Range $target4

# --- Change 5: rewrite the trailing 3 paragraphs into 10 new paragraphs -------------
$target5 = @'
<w:p w14:paraId="22522846" w14:textId="481EFFB4" w:rsidR="00515C3C" w:rsidRDefault="00515C3C" w:rsidP="00D203E1"><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Чем это плохо. </w:t></w:r><w:r w:rsidR="00EE1884"><w:t>Если медальон 10 примет пакет медальона 1, он захочет передавать в момент 40. А это неверно.</w:t></w:r></w:p><w:p><w:r><w:t>Попробуем учесть. Будем добавлять к таймеру 7.</w:t></w:r><w:r><w:t xml:space="preserve"> Таким образом, при множителе 4 мы упустим время передачи: 1*4+7=11 </w:t></w:r><w:r><w:t xml:space="preserve">&gt; 8. </w:t></w:r><w:r><w:t>При множителе 8: 1*8+7 = 15</w:t></w:r><w:r><w:t xml:space="preserve"> &lt; 2*8=16. </w:t></w:r><w:r><w:t>Маловато времени, вообще-то.</w:t></w:r></w:p><w:p><w:r><w:t>При длине пакета</w:t></w:r><w:r><w:t xml:space="preserve"> 7 дребезг был в 2 мс. При 8 – почти исчез.</w:t></w:r></w:p><w:p><w:r><w:t>Итак, длина пакета близка к 8, 9 уже слишком много.</w:t></w:r><w:r><w:t xml:space="preserve"> Перейти, что ли, на другое деление частоты таймером? Перейду на 1/64ю.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:i/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:lang w:val="en-US"/></w:rPr><w:t>Later</w:t></w:r></w:p><w:p><w:r><w:t>Время от момента прерывания на передачу до окончания приема пакета – 8.4 мс.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Перешел на быстрый таймер. Второй сигнал почти не дрожит. Множитель 256, длина пакета 131. </w:t></w:r></w:p><w:p><w:r><w:t>При этом множителе рассчетное расстояние между передачами 16.38 мс. Экспериментальное – 16.9 мс. Очевидно, 500 мкс ушло на какие-то инструкции и округления. Уменьшить множитель нельзя, так как тогда расстояние между передачами сократится до 8 с чем-то мс, а это меньше длины приема пакета. Впрочем, об этом мы уже писали.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Прибрался в файлах немного. Выкинул </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>RX</w:t></w:r><w:r><w:t>_</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>needed</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">Оно всегда </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>needed.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$pStart = Get-ParagraphByText $d "Если медальон 10 примет пакет медальона 1"
$pEnd = $pStart
while ($true) {
    $t = $pEnd.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Нужно") { break }
    $pEnd = $pEnd.Next()
}
$fullRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
Replace-ParagraphXml $fullRange $target5

Write-Output "Foxlocket.docx edit applied."
